# Update SoE tables: collapse multi-line "Intervention" cell values (which
# used a carriage return / line break to separate Part/Cohort descriptions)
# into a single line, joined with "; " instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value = "Part A: Ralmitaront; placebo (monotherapy); Part B: Ralmitaront low or high dose; placebo (add-on to current antipsychotics)"
$ws.Range("G9").Value = "Cohort 1: Ulotaront 10mg; Placebo; Cohort 2: Ulotaront 50mg; Placebo"
$ws.Range("G18").Value = "Cohort 1: Ulotaront (50mg/d to 100mg/d); Placebo; Cohort 2: Ulotaront (25mg/d to 100mg/d); Placebo"
